$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the stray data-only row (old row 5: F=3, G=5.3, H=5.3) that doesn't
# belong to any named plant row in the new layout.
$ws.Rows.Item(5).Delete()

# Drop the second header row (old row 2); the two-row header collapses
# into a single row.
$ws.Rows.Item(2).Delete()

# Clear out whatever remained of the old first header row (content and
# formatting) before writing the new combined header.
$ws.Rows.Item(1).ClearContents()
$ws.Rows.Item(1).ClearFormats()

# Write the new single-row header.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# Match the font used elsewhere for header/label text (Arial 9pt), using a
# dedicated style (distinct from the data-label style) for these unit headers.
$ws.Range("F1:K1").Font.Name = "Arial"
$ws.Range("F1:K1").Font.Size = 9
$ws.Range("F1:K1").ShrinkToFit = $false

# Match the selection left in the saved workbook.
$ws.Range("A4:K4").Select() | Out-Null
